# Update schedule with quizzes.
#
# The quiz column (G) shifts down by one row starting at row 3: what used
# to be blank at G3 becomes "Quiz 1", the existing "Quiz 1" at G5 becomes
# "Quiz 2", and so on through G26 which becomes "Quiz 9" (a brand new
# string). Finally G30 (previously blank) becomes "Quiz 10" (also new).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("G3").Value  = "Quiz 1"
$ws.Range("G5").Value  = "Quiz 2"
$ws.Range("G8").Value  = "Quiz 3"
$ws.Range("G11").Value = "Quiz 4"
$ws.Range("G14").Value = "Quiz 5"
$ws.Range("G17").Value = "Quiz 6"
$ws.Range("G20").Value = "Quiz 7"
$ws.Range("G23").Value = "Quiz 8"
$ws.Range("G26").Value = "Quiz 9"
$ws.Range("G30").Value = "Quiz 10"

# Match the author's final view/selection state.
$ws.Range("G30").Select()
